$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.916.32'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.69%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.813.72'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.50%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.81'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.25%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.01%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4645'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.65%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3707'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.00%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07357'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.15%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8721'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.49%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.47'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.50%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.840.37'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.93%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.349'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.514'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.65%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07054'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.54%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.40'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.43%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.04%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008707'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.10%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9996'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.10%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.71'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.03%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.931.90'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.72%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.317'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.00%  '

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.90%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.036.88'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.95%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.909'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.32%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.62'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.30%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.42'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.99%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.141'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.41%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.314'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.18%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.72'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.48%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08916'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.02%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7563'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.64%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.157'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.99%  '

# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.466'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.52%  '

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.919'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.91%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9993'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.01%  '

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.07%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01958'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.75%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05256'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.21%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.434'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.94%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.923'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.56%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5331'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.48%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.212'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.46%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1664'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.22%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.461'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.15%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4959'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.62%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.32'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.04%  '

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.677'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.19%  '

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9992'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.04%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.20'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.26%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06286'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.40%  '
